# Weekly fruit/hortaliza update: insert a new price record as row 41,
# pushing the existing rows 41-74 down to 42-75.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 41 (shifts 41..74 -> 42..75).
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new weekly record.
$ws.Cells.Item(41, 1).Value = 7
$ws.Cells.Item(41, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(41, 3).Value = "Ñuble"
$ws.Cells.Item(41, 4).Value = 44651
$ws.Cells.Item(41, 5).Value = 16
$ws.Cells.Item(41, 6).Value = 100112022
$ws.Cells.Item(41, 7).Value = "Arveja Verde"
$ws.Cells.Item(41, 8).Value = "Perfection"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 60
$ws.Cells.Item(41, 11).Value = 24000
$ws.Cells.Item(41, 12).Value = 25000
$ws.Cells.Item(41, 13).Value = 24500
$ws.Cells.Item(41, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(41, 15).Value = "Carahue"
$ws.Cells.Item(41, 16).Value = 980
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"
